# Update the "Correspond Handoff Datetime" (col E) and
# "Correspond Handback DateTime" (col H) timestamps on row 2 of the
# "zh-cn" and "de-de" worksheets, reflecting a regenerated handback report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-24 08:54:00"
$wsZhCn.Range("H2").Value = "2016-03-24 08:54:29"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-24 08:54:05"
$wsDeDe.Range("H2").Value = "2016-03-24 08:54:36"
